$d = $word.ActiveDocument

# ---------------------------------------------------------------
# 1. Remove the stray _GoBack bookmark from the "Chatterbox is a AI
#    driven..." paragraph - it moves further down the document to the
#    newly-started "AI can hold a conversation with you" bullet.
# ---------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ---------------------------------------------------------------
# 2. Highlight the "Iteration 3" heading paragraph green - both the
#    paragraph mark (pPr/rPr) and the run itself need the highlight
#    attribute, matching a "select whole paragraph incl. paragraph
#    mark, then apply highlight" Word operation. Because directly
#    setting HighlightColorIndex on a range spanning an existing run
#    only ever touches the run, we temporarily empty the paragraph,
#    highlight the (now run-less) paragraph mark, then restore the
#    text and highlight it too.
# ---------------------------------------------------------------
$iter3 = $d.Paragraphs.Item(11)
$textOnly = $iter3.Range
$textOnly.MoveEnd(1, -1) | Out-Null   # exclude the paragraph mark
$iter3Text = $textOnly.Text
$textOnly.Text = ""                    # paragraph is now run-less

$markRange = $d.Paragraphs.Item(11).Range
$markRange.HighlightColorIndex = 4     # wdBrightGreen -> lands on pPr/rPr

$fillRange = $d.Paragraphs.Item(11).Range
$fillRange.InsertBefore($iter3Text)    # put the text back

$runRange = $d.Paragraphs.Item(11).Range
$runRange.MoveEnd(1, -1) | Out-Null    # exclude the paragraph mark again
$runRange.HighlightColorIndex = 4      # wdBrightGreen -> lands on run rPr

# ---------------------------------------------------------------
# 3. Highlight the "AI can hold a conversation with you" paragraph
#    red using the same two-step trick, then drop the _GoBack
#    bookmark right at the start of its text.
# ---------------------------------------------------------------
$aiPara = $d.Paragraphs.Item(12)
$aiTextOnly = $aiPara.Range
$aiTextOnly.MoveEnd(1, -1) | Out-Null
$aiText = $aiTextOnly.Text
$aiTextOnly.Text = ""

$aiMarkRange = $d.Paragraphs.Item(12).Range
$aiMarkRange.HighlightColorIndex = 6   # wdRed -> lands on pPr/rPr

$aiFillRange = $d.Paragraphs.Item(12).Range
$aiFillRange.InsertBefore($aiText)

$aiRunRange = $d.Paragraphs.Item(12).Range
$aiRunRange.MoveEnd(1, -1) | Out-Null
$aiRunRange.HighlightColorIndex = 6    # wdRed -> lands on run rPr

$aiStart = $d.Paragraphs.Item(12).Range.Start
$bmRange = $d.Range($aiStart, $aiStart)
$d.Bookmarks.Add("_GoBack", $bmRange)
